$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.036957476541118
$ws.Cells.Item(2, 4).Value = 1.054211288482672
$ws.Cells.Item(2, 5).Value = 1.035798253899281
$ws.Cells.Item(2, 6).Value = 1.06083688210989
$ws.Cells.Item(2, 9).Value = 1.046532868676721
$ws.Cells.Item(2, 10).Value = 1.042062963983529
$ws.Cells.Item(2, 11).Value = 1.056955124053742
$ws.Cells.Item(2, 12).Value = 1.038593826904688
$ws.Cells.Item(2, 13).Value = 1.063562576839835
$ws.Cells.Item(2, 14).Value = 1.017903570050256
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038009068662963
$ws.Cells.Item(3, 4).Value = 1.054914599969143
$ws.Cells.Item(3, 5).Value = 1.036695617954739
$ws.Cells.Item(3, 6).Value = 1.061715007593348
$ws.Cells.Item(3, 9).Value = 1.046831781892939
$ws.Cells.Item(3, 10).Value = 1.042758222345012
$ws.Cells.Item(3, 11).Value = 1.057471914828258
$ws.Cells.Item(3, 12).Value = 1.039300426791439
$ws.Cells.Item(3, 13).Value = 1.064255045226615
$ws.Cells.Item(3, 14).Value = 1.018139400982565
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.038689734997415
$ws.Cells.Item(4, 4).Value = 1.055369821495843
$ws.Cells.Item(4, 5).Value = 1.037276797264747
$ws.Cells.Item(4, 6).Value = 1.062283680012021
$ws.Cells.Item(4, 9).Value = 1.047023997630755
$ws.Cells.Item(4, 10).Value = 1.043207745244501
$ws.Cells.Item(4, 11).Value = 1.057805769008389
$ws.Cells.Item(4, 12).Value = 1.039757537355634
$ws.Cells.Item(4, 13).Value = 1.064702944037626
$ws.Cells.Item(4, 14).Value = 1.018291751490072
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.038975938739923
$ws.Cells.Item(5, 4).Value = 1.055561227032463
$ws.Cells.Item(5, 5).Value = 1.037521249991565
$ws.Cells.Item(5, 6).Value = 1.062522860446928
$ws.Cells.Item(5, 9).Value = 1.04710451715803
$ws.Cells.Item(5, 10).Value = 1.043396639078057
$ws.Cells.Item(5, 11).Value = 1.057945990262676
$ws.Cells.Item(5, 12).Value = 1.039949680564929
$ws.Cells.Item(5, 13).Value = 1.064891198003588
$ws.Cells.Item(5, 14).Value = 1.018355740059601
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.039023996632349
$ws.Cells.Item(6, 4).Value = 1.055593366612828
$ws.Cells.Item(6, 5).Value = 1.037562301987018
$ws.Cells.Item(6, 6).Value = 1.06256302634413
$ws.Cells.Item(6, 9).Value = 1.047118019842058
$ws.Cells.Item(6, 10).Value = 1.043428350152308
$ws.Cells.Item(6, 11).Value = 1.057969526320391
$ws.Cells.Item(6, 12).Value = 1.039981940711841
$ws.Cells.Item(6, 13).Value = 1.064922804142788
$ws.Cells.Item(6, 14).Value = 1.018366480520937
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.038693559063097
$ws.Cells.Item(7, 4).Value = 1.055372378948008
$ws.Cells.Item(7, 5).Value = 1.037280063163954
$ws.Cells.Item(7, 6).Value = 1.062286875519736
$ws.Cells.Item(7, 9).Value = 1.047025074668668
$ws.Cells.Item(7, 10).Value = 1.043210269590004
$ws.Cells.Item(7, 11).Value = 1.057807643167574
$ws.Cells.Item(7, 12).Value = 1.039760104885672
$ws.Cells.Item(7, 13).Value = 1.064705459665233
$ws.Cells.Item(7, 14).Value = 1.018292606742895
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037312822270023
$ws.Cells.Item(8, 4).Value = 1.054448947855739
$ws.Cells.Item(8, 5).Value = 1.036101413216717
$ws.Cells.Item(8, 6).Value = 1.061133551452642
$ws.Cells.Item(8, 9).Value = 1.046634136239053
$ws.Cells.Item(8, 10).Value = 1.042298003206302
$ws.Cells.Item(8, 11).Value = 1.057129888071483
$ws.Cells.Item(8, 12).Value = 1.038832647459376
$ws.Cells.Item(8, 13).Value = 1.063796635366512
$ws.Cells.Item(8, 14).Value = 1.017983321396729
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.034881444219985
$ws.Cells.Item(9, 4).Value = 1.052822807163124
$ws.Cells.Item(9, 5).Value = 1.034028530420611
$ws.Cells.Item(9, 6).Value = 1.059104869878649
$ws.Cells.Item(9, 9).Value = 1.045936072348728
$ws.Cells.Item(9, 10).Value = 1.040687767216464
$ws.Cells.Item(9, 11).Value = 1.055931465993954
$ws.Cells.Item(9, 12).Value = 1.037197552826156
$ws.Cells.Item(9, 13).Value = 1.062193874019002
$ws.Cells.Item(9, 14).Value = 1.017436433273638
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033261640449137
$ws.Cells.Item(10, 4).Value = 1.051739497357984
$ws.Cells.Item(10, 5).Value = 1.032649369502656
$ws.Cells.Item(10, 6).Value = 1.057754917145619
$ws.Cells.Item(10, 9).Value = 1.045464548185762
$ws.Cells.Item(10, 10).Value = 1.039612480427784
$ws.Cells.Item(10, 11).Value = 1.05512978730187
$ws.Cells.Item(10, 12).Value = 1.036106982355728
$ws.Cells.Item(10, 13).Value = 1.061124547158391
$ws.Cells.Item(10, 14).Value = 1.017070585659103
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032560511281927
$ws.Cells.Item(11, 4).Value = 1.051270612525423
$ws.Cells.Item(11, 5).Value = 1.032052840031739
$ws.Cells.Item(11, 6).Value = 1.057170979535261
$ws.Cells.Item(11, 9).Value = 1.045258919534829
$ws.Cells.Item(11, 10).Value = 1.039146447690004
$ws.Cells.Item(11, 11).Value = 1.054782013456333
$ws.Cells.Item(11, 12).Value = 1.03563463831095
$ws.Cells.Item(11, 13).Value = 1.060661333184482
$ws.Cells.Item(11, 14).Value = 1.016911874519065
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032300119064093
$ws.Cells.Item(12, 4).Value = 1.0510964783059
$ws.Cells.Item(12, 5).Value = 1.031831361572959
$ws.Cells.Item(12, 6).Value = 1.056954170412553
$ws.Cells.Item(12, 9).Value = 1.045182321550001
$ws.Cells.Item(12, 10).Value = 1.038973278583005
$ws.Cells.Item(12, 11).Value = 1.05465273908423
$ws.Cells.Item(12, 12).Value = 1.035459171061701
$ws.Cells.Item(12, 13).Value = 1.060489247396149
$ws.Cells.Item(12, 14).Value = 1.016852877750821
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032355972345011
$ws.Cells.Item(13, 4).Value = 1.05113382928969
$ws.Cells.Item(13, 5).Value = 1.031878864958828
$ws.Cells.Item(13, 6).Value = 1.057000672556935
$ws.Cells.Item(13, 9).Value = 1.04519876195308
$ws.Cells.Item(13, 10).Value = 1.03901042682689
$ws.Cells.Item(13, 11).Value = 1.054680473194038
$ws.Cells.Item(13, 12).Value = 1.035496810174692
$ws.Cells.Item(13, 13).Value = 1.060526161613075
$ws.Cells.Item(13, 14).Value = 1.01686553476326
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032538986398156
$ws.Cells.Item(14, 4).Value = 1.051256217908184
$ws.Cells.Item(14, 5).Value = 1.032034530527522
$ws.Cells.Item(14, 6).Value = 1.057153056156941
$ws.Cells.Item(14, 9).Value = 1.045252592373799
$ws.Cells.Item(14, 10).Value = 1.039132134777475
$ws.Cells.Item(14, 11).Value = 1.054771329544089
$ws.Cells.Item(14, 12).Value = 1.035620134491129
$ws.Cells.Item(14, 13).Value = 1.060647109075723
$ws.Cells.Item(14, 14).Value = 1.016906998732106
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032651752476219
$ws.Cells.Item(15, 4).Value = 1.051331629648258
$ws.Cells.Item(15, 5).Value = 1.03213045437627
$ws.Cells.Item(15, 6).Value = 1.05724695685209
$ws.Cells.Item(15, 9).Value = 1.045285730149437
$ws.Cells.Item(15, 10).Value = 1.039207114617035
$ws.Cells.Item(15, 11).Value = 1.054827296482952
$ws.Cells.Item(15, 12).Value = 1.035696116346152
$ws.Cells.Item(15, 13).Value = 1.060721625180228
$ws.Cells.Item(15, 14).Value = 1.01693254017623
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033308177473395
$ws.Cells.Item(16, 4).Value = 1.051770619913776
$ws.Cells.Item(16, 5).Value = 1.032688973073119
$ws.Cells.Item(16, 6).Value = 1.057793683909079
$ws.Cells.Item(16, 9).Value = 1.045478164443396
$ws.Cells.Item(16, 10).Value = 1.039643400523947
$ws.Cells.Item(16, 11).Value = 1.055152854448969
$ws.Cells.Item(16, 12).Value = 1.036138327780611
$ws.Cells.Item(16, 13).Value = 1.061155285244271
$ws.Cells.Item(16, 14).Value = 1.017081112559916
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.033720004096101
$ws.Cells.Item(17, 4).Value = 1.052046040030569
$ws.Cells.Item(17, 5).Value = 1.033039493192207
$ws.Cells.Item(17, 6).Value = 1.058136792940099
$ws.Cells.Item(17, 9).Value = 1.045598483929343
$ws.Cells.Item(17, 10).Value = 1.039916957163139
$ws.Cells.Item(17, 11).Value = 1.055356897075077
$ws.Cells.Item(17, 12).Value = 1.036415683503936
$ws.Cells.Item(17, 13).Value = 1.061427258880332
$ws.Cells.Item(17, 14).Value = 1.017174228789806
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033960240362675
$ws.Cells.Item(18, 4).Value = 1.052206706709346
$ws.Cells.Item(18, 5).Value = 1.033244008934985
$ws.Cells.Item(18, 6).Value = 1.058336980534973
$ws.Cells.Item(18, 9).Value = 1.045668523785215
$ws.Cells.Item(18, 10).Value = 1.040076476956734
$ws.Cells.Item(18, 11).Value = 1.055475849725275
$ws.Cells.Item(18, 12).Value = 1.036577448818128
$ws.Cells.Item(18, 13).Value = 1.061585878270007
$ws.Cells.Item(18, 14).Value = 1.01722851327027
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034042158927193
$ws.Cells.Item(19, 4).Value = 1.052261493046897
$ws.Cells.Item(19, 5).Value = 1.033313754236237
$ws.Cells.Item(19, 6).Value = 1.058405249121115
$ws.Cells.Item(19, 9).Value = 1.045692381730643
$ws.Cells.Item(19, 10).Value = 1.040130862093262
$ws.Cells.Item(19, 11).Value = 1.055516398973193
$ws.Cells.Item(19, 12).Value = 1.03663260463351
$ws.Cells.Item(19, 13).Value = 1.061639960265903
$ws.Cells.Item(19, 14).Value = 1.017247018011388
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033675816436528
$ws.Cells.Item(20, 4).Value = 1.052016488107514
$ws.Cells.Item(20, 5).Value = 1.033001879150578
$ws.Cells.Item(20, 6).Value = 1.058099974591601
$ws.Cells.Item(20, 9).Value = 1.045585589313651
$ws.Cells.Item(20, 10).Value = 1.039887611380599
$ws.Cells.Item(20, 11).Value = 1.055335011636034
$ws.Cells.Item(20, 12).Value = 1.036385927055469
$ws.Cells.Item(20, 13).Value = 1.061398080562921
$ws.Cells.Item(20, 14).Value = 1.017164241260642
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032485092263095
$ws.Cells.Item(21, 4).Value = 1.051220176660008
$ws.Cells.Item(21, 5).Value = 1.031988688156081
$ws.Cells.Item(21, 6).Value = 1.057108180455994
$ws.Cells.Item(21, 9).Value = 1.045236746681363
$ws.Cells.Item(21, 10).Value = 1.039096296578457
$ws.Cells.Item(21, 11).Value = 1.054744577248993
$ws.Cells.Item(21, 12).Value = 1.035583819037736
$ws.Cells.Item(21, 13).Value = 1.060611493814979
$ws.Cells.Item(21, 14).Value = 1.016894789850415
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031736658646209
$ws.Cells.Item(22, 4).Value = 1.050719680972363
$ws.Cells.Item(22, 5).Value = 1.031352228669525
$ws.Cells.Item(22, 6).Value = 1.056485128897511
$ws.Cells.Item(22, 9).Value = 1.045016152344466
$ws.Cells.Item(22, 10).Value = 1.038598396575564
$ws.Cells.Item(22, 11).Value = 1.054372794491405
$ws.Cells.Item(22, 12).Value = 1.035079400013419
$ws.Cells.Item(22, 13).Value = 1.06011677683637
$ws.Cells.Item(22, 14).Value = 1.016725118356522
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032133396487105
$ws.Cells.Item(23, 4).Value = 1.050984986062921
$ws.Cells.Item(23, 5).Value = 1.031689573304054
$ws.Cells.Item(23, 6).Value = 1.056815369848983
$ws.Cells.Item(23, 9).Value = 1.04513321319378
$ws.Cells.Item(23, 10).Value = 1.038862377671096
$ws.Cells.Item(23, 11).Value = 1.054569935763231
$ws.Cells.Item(23, 12).Value = 1.03534681161721
$ws.Cells.Item(23, 13).Value = 1.060379050437226
$ws.Cells.Item(23, 14).Value = 1.016815088693265
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.033695782861168
$ws.Cells.Item(24, 4).Value = 1.052029841291548
$ws.Cells.Item(24, 5).Value = 1.033018875121346
$ws.Cells.Item(24, 6).Value = 1.058116611039619
$ws.Cells.Item(24, 9).Value = 1.045591416269445
$ws.Cells.Item(24, 10).Value = 1.039900871603692
$ws.Cells.Item(24, 11).Value = 1.055344900915381
$ws.Cells.Item(24, 12).Value = 1.036399372749442
$ws.Cells.Item(24, 13).Value = 1.061411265044961
$ws.Cells.Item(24, 14).Value = 1.017168754283654
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035509817307851
$ws.Cells.Item(25, 4).Value = 1.053243070395033
$ws.Cells.Item(25, 5).Value = 1.034563936503253
$ws.Cells.Item(25, 6).Value = 1.059628896578428
$ws.Cells.Item(25, 9).Value = 1.046117624038343
$ws.Cells.Item(25, 10).Value = 1.041104370104771
$ws.Cells.Item(25, 11).Value = 1.056241771544499
$ws.Cells.Item(25, 12).Value = 1.037620355456714
$ws.Cells.Item(25, 13).Value = 1.062608374499816
$ws.Cells.Item(25, 14).Value = 1.017578039264079
